$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.004.02"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.371.16"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.41"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.77%  "
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.67"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "37.19"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +14.84%  "
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.108"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.40"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.930"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.69%  "
$ws.Range("D17").Value = "2.377.37"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "43.982.19"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  +2.07%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.63"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "78.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "254.09"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  +3.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.89"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.44"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.98"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  +1.56%  "
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0761"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.44"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.64"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.75%  "
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("E39").Value = "  +2.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.52"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +17.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.77"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +9.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.86"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +14.16%  "
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.205"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.10"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("E46").Value = "  +3.79%  "
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.16"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.76"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.42"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +17.06%  "
